$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 29 (record #27) ---
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = 26243521
$ws.Range("C29").Value = 111925
$ws.Range("D29").Value = 105999.85
$ws.Range("E29").Value = 44237
$ws.Range("F29").Formula = "=IF(B29=`"`",`"`",C29-D29)"
$ws.Range("G29").Formula = "=IF(B29=`"`",`"`",F29/D29*100)"
$ws.Range("H29").Formula = "=IF(B29=`"`",`"`",D29*1.04)"
$ws.Range("I29").Formula = "=IF(B29=`"`",`"`",C29-H29)"

# --- Row 30 (record #28) ---
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = 26268624
$ws.Range("C30").Value = 32733
$ws.Range("D30").Value = 31000.51
$ws.Range("E30").Value = 44241
$ws.Range("F30").Formula = "=IF(B30=`"`",`"`",C30-D30)"
$ws.Range("G30").Formula = "=IF(B30=`"`",`"`",F30/D30*100)"
$ws.Range("H30").Formula = "=IF(B30=`"`",`"`",D30*1.04)"
$ws.Range("I30").Formula = "=IF(B30=`"`",`"`",C30-H30)"

# Match the formatting of the preceding data row (28) for the A/C/D/E
# columns (B already inherits the correct column default style, and F:I
# already carried the right style since they pre-existed as blank shared
# formula cells).
$ws.Range("A28").Copy()
$ws.Range("A29:A30").PasteSpecial(-4122)
$ws.Range("C28").Copy()
$ws.Range("C29:C30").PasteSpecial(-4122)
$ws.Range("D28").Copy()
$ws.Range("D29:D30").PasteSpecial(-4122)
$ws.Range("E28").Copy()
$ws.Range("E29:E30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the frozen-pane/selection view state ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A13").Select()
$ws.Range("I33").Select()
